$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'26.912.02"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.Value = "'  -0.29%  "
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.Value = "'1.549.90"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.Value = "'  -0.38%  "
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.Value = "'  -0.18%  "
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.Value = "'206.54"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.Value = "'  -0.17%  "
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.Value = "'0.486"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.Value = "'  +0.06%  "
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.Value = "'  -0.21%  "
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.Value = "'22.10"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.Value = "'  +2.64%  "
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.Value = "'0.245"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.Value = "'  -0.85%  "
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.Value = "'  +0.46%  "
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.Value = "'0.0855"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.Value = "'  -0.45%  "
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.Value = "'  -0.32%  "
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.Value = "'1.552.23"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.Value = "'  -0.08%  "
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.Value = "'  +0.63%  "
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.Value = "'0.518"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.Value = "'  +0.74%  "
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.Value = "'26.898.88"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.Value = "'  -0.29%  "
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.Value = "'61.63"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.Value = "'  -0.08%  "
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.Value = "'217.30"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.Value = "'  +1.17%  "
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.Value = "'  +1.27%  "
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.Value = "'7.24"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.Value = "'  +0.02%  "
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.Value = "'  -0.15%  "
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.Value = "'  +0.08%  "
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.Value = "'  -0.16%  "
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.Value = "'  -0.77%  "
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.Value = "'154.21"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.Value = "'  +0.24%  "
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.Value = "'  -0.71%  "
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.Value = "'  +0.35%  "
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.Value = "'  +0.69%  "
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.Value = "'  -0.24%  "
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.Value = "'  +1.21%  "
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.Value = "'  -1.16%  "
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.Value = "'3.22"
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.Value = "'  -0.16%  "
$cell.ClearFormats()
$cell = $ws.Range("D33")
$cell.Value = "'1.414.69"
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.Value = "'  +2.91%  "
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.Value = "'  +4.19%  "
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.Value = "'  +2.31%  "
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.Value = "'0.972"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.Value = "'  +0.11%  "
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.Value = "'  +0.15%  "
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.Value = "'  +0.25%  "
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.Value = "'0.522"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.Value = "'  +0.39%  "
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.Value = "'  -0.27%  "
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.Value = "'  +4.80%  "
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.Value = "'  -0.29%  "
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.Value = "'0.996"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.Value = "'  +1.56%  "
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.Value = "'64.26"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.Value = "'  +0.51%  "
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.Value = "'1.75"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.Value = "'  +0.03%  "
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.Value = "'1.684.39"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.Value = "'  -0.33%  "
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.Value = "'87.53"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.Value = "'  +1.45%  "
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.Value = "'  +2.04%  "
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.Value = "'  +8.34%  "
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.Value = "'0.0954"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.Value = "'  -0.18%  "
$cell.ClearFormats()
